{"js": "// Insert a new \"Replace hard coded ids ...\" list item right after the\n// \"Expense form is blinking in Chrome and Explorer\" paragraph. Using\n// Paragraph.insertParagraph(text, \"After\") copies the anchor paragraph's\n// formatting (ListParagraph style, level-1 numbered list numId 3,\n// justified, HTMLCode run style) onto the new paragraph/run automatically.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"Expense form is blinking in Chrome and Explorer\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === targetText) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find anchor paragraph: \" + targetText);\n}\n\nconst newText = \"Replace hard coded ids (such as 1) with fragment of code in all tests\";\ntarget.insertParagraph(newText, \"After\");\n\nawait context.sync();\n", "ps1": "# Insert a new \"Replace hard coded ids ...\" list item right after the\n# \"Expense form is blinking in Chrome and Explorer\" paragraph, inheriting\n# that paragraph's list/formatting (ListParagraph style, level-1 numbered\n# list numId 3, justified, HTMLCode run style).\n\n$d = $word.ActiveDocument\n\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.TrimEnd(\"`r`a\") -eq \"Expense form is blinking in Chrome and Explorer\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find anchor paragraph\"\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n$target.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($targetIndex + 1)\n$newPara.Range.Text = \"Replace hard coded ids (such as 1) with fragment of code in all tests\"\n"}
